# Reorders the docente columns C:J (and their column widths) according to
# the new layout, and fixes a small content typo in row 7 (order of the
# two "BATISTA, DIEGO ..." entries in the merged citation list).
#
# Mapping of OLD column letter -> NEW column letter (where that column's
# data ends up):
#   C -> I    D -> F    E -> C    F -> H
#   G -> J    H -> E    I -> G    J -> D

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Old -> New column index map (1-based). Columns A=1, B=2, C=3, ... J=10
$colMap = @{
    3  = 9   # C -> I
    4  = 6   # D -> F
    5  = 3   # E -> C
    6  = 8   # F -> H
    7  = 10  # G -> J
    8  = 5   # H -> E
    9  = 7   # I -> G
    10 = 4   # J -> D
}

# Snapshot all current values (rows 1..lastRow, cols C..J) before writing
# anything, so that writes to the new positions don't clobber values we
# still need to read for later columns.
$snapshot = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    foreach ($oldCol in $colMap.Keys) {
        $snapshot["$r,$oldCol"] = $ws.Cells.Item($r, $oldCol).Value2
    }
}

# Write the values into their new homes.
for ($r = 1; $r -le $lastRow; $r++) {
    foreach ($oldCol in $colMap.Keys) {
        $newCol = $colMap[$oldCol]
        $val = $snapshot["$r,$oldCol"]
        if ($null -eq $val) { $val = "" }
        $ws.Cells.Item($r, $newCol).Value = $val
    }
}

# Fix the citation text for Diego Silva Batista (row 7, now column C):
# the two alternate-name entries were swapped in order.
$ws.Cells.Item(7, 3).Value = "BATISTA, D. S. (1) | BATISTA, DIEGO SILVA (1) | BATISTA, DIEGO S (1)"

# Re-apply the column widths so they travel with the moved content.
# Note: Excel's ColumnWidth (COM, character units) differs from the
# OOXML <col width="..."> attribute by a constant ~0.83 offset caused
# by cell padding, so we subtract that before assigning.
$widths = @{
    3  = 70
    4  = 24
    5  = 21
    6  = 27
    7  = 38
    8  = 28
    9  = 28
    10 = 25
}
foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col] - 0.83
}

$wb.Save()
